$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns into the result table -----------------------
# 1) A brand-new column A (will hold "Date"), pushing the old A..L to B..M
$ws.Columns.Item(1).Insert()
# 2) A brand-new column D (will hold "loss"), pushing the (now shifted)
#    D..M columns to E..N
$ws.Columns.Item(4).Insert()

# --- New header row values ------------------------------------------------
# Set D1 ("loss") first so the shared-string table fills up in the same
# order the source workbook uses (loss, Date, Run-Name ...).
$ws.Range("D1").Value = "loss"
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Run-Name (for the Log-file)"

# --- Row 3 ("first run (no changes)") gets split across two columns ------
$ws.Range("B3").Value = "first run"
$ws.Range("C3").Value = "no changes"

# --- Row 4 gets a real date value in the new column A ---------------------
$ws.Range("A4").Value = 43851
$ws.Range("A4").NumberFormat = "mm-dd-yy"

# --- Column widths (best effort match of the resized/auto-fit columns) ---
$ws.Columns.Item(1).ColumnWidth = 8.75
$ws.Columns.Item(2).ColumnWidth = 21.25
$ws.Columns.Item(4).ColumnWidth = 10.75

# --- Selection moves to A5 -------------------------------------------------
$ws.Range("A5").Select() | Out-Null
